# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit columns (H:N) for the affected Leve rows on each
# job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Labels/IDs in A:G are
# untouched; only the numeric market snapshot columns move.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 3740
$ws.Range("I8").Value = 3740
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 11220
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -11081
$ws.Range("N8").Value = $null

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 7518.3335
$ws.Range("I29").Value = 5722
$ws.Range("J29").Value = 11111
$ws.Range("K29").Value = 17166
$ws.Range("L29").Value = 33333
$ws.Range("M29").Value = -16885
$ws.Range("N29").Value = -33895

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 737
$ws.Range("J38").Value = 1368.1666
$ws.Range("L38").Value = 4104.4998
$ws.Range("N38").Value = -4848.4998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 403.5
$ws.Range("I58").Value = 204.66667
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 614.00001
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -464.00001
$ws.Range("N58").Value = -3300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 11541187
$ws.Range("I116").Value = 23072824
$ws.Range("J116").Value = 9549.667
$ws.Range("K116").Value = 23072824
$ws.Range("L116").Value = 9549.667
$ws.Range("M116").Value = -23069382
$ws.Range("N116").Value = -16433.667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 298351.6
$ws.Range("I132").Value = 338809.28
$ws.Range("J132").Value = 75834.336
$ws.Range("K132").Value = 1016427.84
$ws.Range("L132").Value = 227503.008
$ws.Range("M132").Value = -1013897.84
$ws.Range("N132").Value = -232563.008

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 17333.334
$ws.Range("J133").Value = 17333.334
$ws.Range("L133").Value = 17333.334
$ws.Range("N133").Value = -27453.334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 35010.668
$ws.Range("I21").Value = 45007.5
$ws.Range("J21").Value = 15017
$ws.Range("K21").Value = 45007.5
$ws.Range("L21").Value = 15017
$ws.Range("M21").Value = -44633.5
$ws.Range("N21").Value = -15765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1016
$ws.Range("I35").Value = 1016
$ws.Range("K35").Value = 1016
$ws.Range("M35").Value = -610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1407.2222
$ws.Range("I45").Value = 917.6923
$ws.Range("J45").Value = 2680
$ws.Range("K45").Value = 917.6923
$ws.Range("L45").Value = 2680
$ws.Range("M45").Value = -540.6923
$ws.Range("N45").Value = -3434

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1332.04
$ws.Range("I122").Value = 1435.55
$ws.Range("J122").Value = 918
$ws.Range("K122").Value = 4306.65
$ws.Range("L122").Value = 2754
$ws.Range("M122").Value = -1856.65
$ws.Range("N122").Value = -7654

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1527.5
$ws.Range("I99").Value = 1527.5
$ws.Range("K99").Value = 1527.5
$ws.Range("M99").Value = -29.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 41440
$ws.Range("J124").Value = 41440
$ws.Range("L124").Value = 41440
$ws.Range("N124").Value = -51260

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 83685
$ws.Range("I16").Value = 100322.2
$ws.Range("J16").Value = 499
$ws.Range("K16").Value = 100322.2
$ws.Range("L16").Value = 499
$ws.Range("M16").Value = -100035.2
$ws.Range("N16").Value = -1073

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 41633.332
$ws.Range("J59").Value = 41633.332
$ws.Range("L59").Value = 41633.332
$ws.Range("N59").Value = -43923.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10418288
$ws.Range("I99").Value = 12501562
$ws.Range("J99").Value = 1914
$ws.Range("K99").Value = 12501562
$ws.Range("L99").Value = 1914
$ws.Range("M99").Value = -12500064
$ws.Range("N99").Value = -4910

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 610.2857
$ws.Range("I105").Value = 503.46155
$ws.Range("K105").Value = 503.46155
$ws.Range("M105").Value = 1243.53845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 83685
$ws.Range("I113").Value = 100322.2
$ws.Range("J113").Value = 499
$ws.Range("K113").Value = 100322.2
$ws.Range("L113").Value = 499
$ws.Range("M113").Value = -98152.2
$ws.Range("N113").Value = -4839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10418288
$ws.Range("I126").Value = 12501562
$ws.Range("J126").Value = 1914
$ws.Range("K126").Value = 37504686
$ws.Range("L126").Value = 5742
$ws.Range("M126").Value = -37502216
$ws.Range("N126").Value = -10682

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 1600
$ws.Range("J124").Value = 1600
$ws.Range("L124").Value = 4800
$ws.Range("N124").Value = -14620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1822.6586
$ws.Range("I131").Value = 567.8
$ws.Range("J131").Value = 1996.9445
$ws.Range("K131").Value = 1703.4
$ws.Range("L131").Value = 5990.833500000001
$ws.Range("M131").Value = 3336.6
$ws.Range("N131").Value = -16070.8335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1329.4117
$ws.Range("I132").Value = 866.6667
$ws.Range("J132").Value = 1428.5714
$ws.Range("K132").Value = 7800.0003
$ws.Range("L132").Value = 12857.1426
$ws.Range("M132").Value = -5270.0003
$ws.Range("N132").Value = -17917.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10520.929
$ws.Range("I99").Value = 9791.77
$ws.Range("K99").Value = 9791.77
$ws.Range("M99").Value = -7545.77

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1140.2
$ws.Range("I102").Value = 799.7143
$ws.Range("J102").Value = 1934.6666
$ws.Range("K102").Value = 799.7143
$ws.Range("L102").Value = 1934.6666
$ws.Range("M102").Value = 822.2857
$ws.Range("N102").Value = -5178.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2266.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 927140.06
$ws.Range("I122").Value = 2222734
$ws.Range("J122").Value = 1715.8572
$ws.Range("K122").Value = 6668202
$ws.Range("L122").Value = 5147.571599999999
$ws.Range("M122").Value = -6665752
$ws.Range("N122").Value = -10047.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 918.73914
$ws.Range("I16").Value = 936.45
$ws.Range("K16").Value = 936.45
$ws.Range("M16").Value = -766.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9353
$ws.Range("I22").Value = 780.4
$ws.Range("J22").Value = 15476.286
$ws.Range("K22").Value = 780.4
$ws.Range("L22").Value = 15476.286
$ws.Range("M22").Value = -485.4
$ws.Range("N22").Value = -16066.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 9353
$ws.Range("I27").Value = 780.4
$ws.Range("J27").Value = 15476.286
$ws.Range("K27").Value = 780.4
$ws.Range("L27").Value = 15476.286
$ws.Range("M27").Value = -673.4
$ws.Range("N27").Value = -15690.286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 478.2857
$ws.Range("I107").Value = 596.6667
$ws.Range("J107").Value = 389.5
$ws.Range("K107").Value = 1790.0001
$ws.Range("L107").Value = 1168.5
$ws.Range("M107").Value = 129.9999
$ws.Range("N107").Value = -5008.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 231.66667
$ws.Range("I113").Value = 231.66667
$ws.Range("K113").Value = 695.00001
$ws.Range("M113").Value = 1474.99999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 45613.26
$ws.Range("I122").Value = 126288
$ws.Range("J122").Value = 2586.7334
$ws.Range("K122").Value = 378864
$ws.Range("L122").Value = 7760.2002
$ws.Range("M122").Value = -376414
$ws.Range("N122").Value = -12660.2002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10640796
$ws.Range("I132").Value = 16130815
$ws.Range("K132").Value = 48392445
$ws.Range("M132").Value = -48389915
